$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (shifts old rows 35-64 down to 36-65)
$ws.Rows.Item(35).Insert()

# Insert another new row at position 65 (shifts old row 64, now at 65, down to 66)
$ws.Rows.Item(65).Insert()

# Row 35
$ws.Range("A35").NumberFormat = "@"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("A35").Value = "-406"
$ws.Range("B35").Value = "5/8/2025"
$ws.Range("C35").Value = "Olof palme 4144"
$ws.Range("D35").Value = "12"
$ws.Range("E35").Value = "805791925"
$ws.Range("F35").Value = "NEW"
$ws.Range("G35").Value = "Pendiente"
$ws.Range("H35").Value = "Tensar 2 riendas a pique columna 168"
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = "Tensor"
$ws.Range("K35").Value = "Sin equipos"
$ws.Range("L35").Value = "Terminal"
$ws.Range("M35").Value = -58.488252
$ws.Range("N35").Value = -34.553391
$ws.Range("O35").Value = "Saavedra"
$ws.Range("P35").Value = "Capital Norte"

# Row 65
$ws.Range("A65").NumberFormat = "@"
$ws.Range("B65").NumberFormat = "@"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("E65").NumberFormat = "@"
$ws.Range("A65").Value = "-536"
$ws.Range("B65").Value = "7/29/2025"
$ws.Range("C65").Value = "Olof palme 4142"
$ws.Range("D65").Value = "12"
$ws.Range("E65").Value = "ICD30249764 "
$ws.Range("F65").Value = "NEW"
$ws.Range("G65").Value = "Pendiente"
$ws.Range("H65").Value = "Aplomar o desmontar poste"
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = "Desmonte"
$ws.Range("K65").Value = "Sin equipos"
$ws.Range("L65").Value = "Poste"
$ws.Range("M65").Value = -58.488239
$ws.Range("N65").Value = -34.55341
$ws.Range("O65").Value = "Saavedra"
$ws.Range("P65").Value = "Capital Norte"

# Row 67
$ws.Range("A67").NumberFormat = "@"
$ws.Range("B67").NumberFormat = "@"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("E67").NumberFormat = "@"
$ws.Range("A67").Value = "-543"
$ws.Range("B67").Value = "8/1/2025"
$ws.Range("C67").Value = "Pedro Ignacio Rivera 3258"
$ws.Range("D67").Value = "13"
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = "NEW"
$ws.Range("G67").Value = "Pendiente"
$ws.Range("H67").Value = "Desmontar poste en desuso"
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = "Desmonte"
$ws.Range("K67").Value = "Sin equipos"
$ws.Range("L67").Value = "Poste"
$ws.Range("M67").Value = -58.46967
$ws.Range("N67").Value = -34.561676
$ws.Range("O67").Value = "Colegiales"
$ws.Range("P67").Value = "Capital Norte"
